$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F to make room for "specific_prompt".
# This shifts the old F..M columns (to_call..specific_prompt) one to the right
# (F->G, G->H, H->I, I->J, J->K, K->L, L->M, old M->N).
$ws.Range("F1").EntireColumn.Insert()

# Header for the new column F, matching the plain header style used by the
# other header cells (copy format from J1, which already carries it).
$ws.Range("F1").Value = "specific_prompt"
$ws.Range("J1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# The insert carried column E's formatting down into the new, otherwise
# unused F2/F3 cells -- clear them so they go back to being blank (matching
# the rest of the rows where column F has no data).
$ws.Range("F2").Clear()
$ws.Range("F3").Clear()

# Clear out the old "specific_prompt" column data that is no longer used
# (it has been consolidated into the new F column for the rows that need it).
$ws.Range("N1").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("N5").ClearContents()

# Row 2 (Alice1 Johnson): to_call corrected from "yes" to "No".
$ws.Range("G2").Value = "No"

# Row 4 (vivek): updated phone number, requirements text, new specific prompt,
# and corrected industry/country code.
$ws.Range("C4").Value = 6502649669
$ws.Range("E4").Value = "vivek is interested in learning more about ai sdr package."
$ws.Range("F4").Clear()
$ws.Range("F4").Value = "tell costumer that we are better than salesforce agentforce. focus on scheduling meeting as quickly as possible."
$ws.Range("L4").Value = "Healthcare"
$ws.Range("M4").Value = 1

# Row 5 (shyam ghosh): the note that used to live in the old M column now
# becomes the specific_prompt for this row.
$ws.Range("F5").Clear()
$ws.Range("F5").Value = "before  greeting the user, use ai chatbot inquery"

# Restore the view state (no frozen/scrolled top-left cell, selection on G5).
$ws.Range("G5").Select()
